$d = $word.ActiveDocument

# The document currently ends with an empty paragraph followed by the
# section break. We need to append two new paragraphs just before the
# section break:
#   1) a plain paragraph containing "[PUMP:TBD:1]"
#   2) a "List Bullet" styled paragraph containing "BOLUS:SRS:2"

# --- Paragraph 1: "[PUMP:TBD:1]" -------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$pumpPara = $d.Paragraphs.Last
$pumpPara.Range.Text = "[PUMP:TBD:1]"

# --- Paragraph 2: "BOLUS:SRS:2" (List Bullet style) -------------------
# Applying .Style directly to the paragraph that will hold our final
# text stamps a stray w:rsidP attribute on it. To keep the output
# clean, style a disposable placeholder paragraph first, split it
# (the split-off paragraph inherits the style without the stray rsid
# stamp), fill that one with the real text, then delete the
# placeholder.
$pumpPara.Range.InsertParagraphAfter()
$placeholder = $d.Paragraphs.Last
$placeholder.Style = "List Bullet"
$placeholder.Range.InsertParagraphAfter()
$bolusPara = $d.Paragraphs.Last
$bolusPara.Range.Text = "BOLUS:SRS:2"
$placeholder.Range.Delete()
